# Methods-diagram poster: bump the cached "datetimeFigureOut" field text
# (07/04/2014 -> 08/04/2014) on the slide master and every slide layout,
# and split the "3 CNV events" run into "3 " / "CNA " / "events" on the
# slide itself.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, [string]$newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$design = $p.Designs.Item(1)
$slideMaster = $design.SlideMaster

# Slide master's own Date Placeholder.
Update-DatePlaceholder $slideMaster.Shapes "08/04/2014"

# Every custom (slide) layout's Date Placeholder.
$layouts = $slideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "08/04/2014"
}

# Slide 1: split "3 CNV events" into three runs: "3 " / "CNA " / "events".
$slide = $p.Slides.Item(1)
$target = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 28") {
        $target = $sh
    } elseif ($target -eq $null -and $sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -like "*CNV events*") {
        $target = $sh
    }
}

if ($target -ne $null) {
    $tr = $target.TextFrame.TextRange
    $fullText = $tr.Text
    $firstLine = $fullText -replace "3 CNV events$", ""
    # First drop the "CNV " word (intermediate commit) so the replacement
    # below is recognised as a fresh word insertion rather than a
    # character-level diff against "CNV" - this reproduces the exact
    # "3 " / "CNA " / "events" run split seen in the authored edit.
    $tr.Text = $firstLine + "3 events"
    $tr.Text = $firstLine + "3 CNA events"
}
